$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to be treated as text so date-like strings
    # (e.g. "2022-06-02") are not auto-converted into date serials.
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $value
    $ws.Range($cell).ClearFormats()
}

# Remove the old C1 numeric value - the new layout only uses columns A and B.
$ws.Range("C1").ClearContents()

# Row 1
Set-TextValue "A1" "2022-06-02"
Set-TextValue "B1" "<p>Analisa tanggal 2 juni 2</p>"

# Row 2
Set-TextValue "A2" "2022-06-01"
Set-TextValue "B2" "<p>ASSSSSSSSSSSSSS</p>"

# Row 3
Set-TextValue "A3" "2022-05-31"
Set-TextValue "B3" "<p>Naiknya harga BBM dsa</p>"

# Row 4
Set-TextValue "A4" "2022-05-31"
Set-TextValue "B4" "<p>Analisa tanggal 2 juni</p>"

# Row 5
Set-TextValue "A5" "2022-05-30"
Set-TextValue "B5" "<p>Ww</p>"

# Row 6
Set-TextValue "A6" "2022-05-30"
Set-TextValue "B6" "<p>Re</p>"

# Row 7
Set-TextValue "A7" "2022-05-30"
Set-TextValue "B7" "<p>Lumen (8.3.4) (Laravel Components ^8.0)</p>"
